$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.334.71"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +2.46%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.691.86"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.20%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "524.13"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.92%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "147.13"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.10%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.997"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.20%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.576"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.83%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.715.80"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.86%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "6.50"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.48%  "

$ws.Range("E11").Value = "  -0.07%  "

$ws.Range("E12").Value = "  +1.39%  "

$ws.Range("E13").Value = "  +1.56%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.168.32"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.20%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "60.370.98"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.52%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.838.60"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +5.67%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "21.41"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +2.00%  "

$ws.Range("E18").Value = "  +1.19%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "351.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +1.52%  "

$ws.Range("E20").Value = "  +0.30%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.59"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.36"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.37%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.997"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.24%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.04"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +3.46%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.423"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.86%  "

$ws.Range("E26").Value = "  +5.33%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.994"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.39%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.0₃0821"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.29%  "

$ws.Range("E29").Value = "  +2.15%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.87"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.62%  "

$ws.Range("E31").Value = "  +0.15%  "

$ws.Range("E32").Value = "  +1.75%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "19.17"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.07%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "147.60"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.68%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.30"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +6.61%  "

$ws.Range("E36").Value = "  +9.93%  "

$ws.Range("E37").Value = "  -5.23%  "

$ws.Range("E38").Value = "  +9.26%  "

$ws.Range("E39").Value = "  +4.16%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "36.93"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.43%  "

$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "286.34"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.21%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "20.14"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +1.63%  "

$ws.Range("E44").Value = "  -0.65%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0992"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.13%  "

$ws.Range("E46").Value = "  +0.30%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.138.86"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.82%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0541"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.49%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "4.91"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +3.92%  "

$ws.Range("E50").Value = "  +2.26%  "

$ws.Range("B51").Value = "WhiteBITCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "10.45"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.77%  "
